$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Sistema de controle de Cópias" paragraph: merge the three runs (plain
#    text / underlined "cópias" / plain text) into a single plain run with
#    the full sentence (drop the underline formatting on "cópias").
#    Replacing the exact same text via Find/Replace with the first run's
#    formatting collapses the run boundaries into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Sistema que pode fornecer informações de quantas cópias foram realizadas em um dia, em uma semana, no mês ou até mesmo no ano. Podendo fazer um levantamento de quanto de papel A4 ou outro tipo foi utilizado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sistema que pode fornecer informações de quantas cópias foram realizadas em um dia, em uma semana, no mês ou até mesmo no ano. Podendo fazer um levantamento de quanto de papel A4 ou outro tipo foi utilizado.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Pedidos via Redes Sociais
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Sistema que permite que pedidos sejam feitos via redes sociais, o sistema em ligação com as redes sociais armazena e efetua o pedido.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sistema que permite que pedidos sejam feitos via redes sociais, armazenando e efetuando o pedido.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Marketing para produtos
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Divulgação do trabalho da empresa, para que tenha um lucro maior por parte da empresa, e esse marketing pode ser feito pelo site da empresa, Anúncios em redes sociais.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Divulgação do trabalho da empresa, para que tenha um lucro maior por parte da empresa, e esse marketing pode ser feito pelo site da empresa, ou anúncios em redes sociais.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Anúncios em redes sociais (1st occurrence)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Divulgar os produtos via redes sociais para maior alcance de pessoas e clientes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Divulgar os produtos e a gráfica via redes sociais tendo em mente um número maior de clientes.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Promoção via E-mail: text changes and the trailing tab is removed
#    (leading tab before the sentence remains). Replace the sentence text
#    first, then delete the single trailing tab character that follows it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Envio de mala direta para clientes para que possam receber promoções.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Envio de promoções para clientes que realizaram um cadastro no site. ",
    2) | Out-Null

$rTrailingTab = $d.Content
$foundTrailingTab = $rTrailingTab.Find.Execute(
    "Envio de promoções para clientes que realizaram um cadastro no site. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundTrailingTab) {
    $tabPos = $rTrailingTab.End
    $tabRange = $d.Range($tabPos, $tabPos + 1)
    if ($tabRange.Text -eq "`t") {
        $tabRange.Delete()
    }
}

# ---------------------------------------------------------------------------
# 6) Expansão Internacional: text changes, and a new empty paragraph is
#    inserted right after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Abertura de uma franquia da empresa em outro país, para expansão do negócio.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Abertura de uma franquia da empresa em outro país, visando uma expansão do negócio e um lucro maior.",
    2) | Out-Null

$pPrFirstLine708 = "<w:pPr><w:spacing w:before='0' w:after='160' w:line='259'/><w:ind w:right='0' w:left='0' w:firstLine='708'/><w:jc w:val='both'/><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial' w:eastAsia='Arial'/><w:color w:val='auto'/><w:spacing w:val='0'/><w:position w:val='0'/><w:sz w:val='24'/><w:shd w:fill='auto' w:val='clear'/></w:rPr></w:pPr>"
$emptyParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>$pPrFirstLine708</w:p>"

$rFranquia = $d.Content
$foundFranquia = $rFranquia.Find.Execute(
    "Abertura de uma franquia da empresa em outro país, visando uma expansão do negócio e um lucro maior.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundFranquia) {
    $insertPoint = $rFranquia.End
    $newRange = $d.Range($insertPoint, $insertPoint)
    $newRange.InsertXML($emptyParaXml) | Out-Null
}

# ---------------------------------------------------------------------------
# 7) Impressão Delivery
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Entregas de encomendas e impressão com pedidos online.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Entregas de impressões, solicitadas online.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Office Boy para entregas
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Entregas feitas por funcionário em horário comercial.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Entregas de pedidos feitos pelo website, realizadas por funcionário em horário comercial.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 9) Enquete de Satisfação do Cliente: text changes, and TWO new empty
#    paragraphs are inserted right after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Enquete para saber se os clientes estão satisfeitos com o produto entregue.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enquete para saber se os clientes estão satisfeitos com o produto e o serviço que foi prestado.",
    2) | Out-Null

$rEnquete = $d.Content
$foundEnquete = $rEnquete.Find.Execute(
    "Enquete para saber se os clientes estão satisfeitos com o produto e o serviço que foi prestado.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundEnquete) {
    $insertPoint1 = $rEnquete.End
    $newRange1 = $d.Range($insertPoint1, $insertPoint1)
    $newRange1.InsertXML($emptyParaXml) | Out-Null
}

$rEnquete2 = $d.Content
$foundEnquete2 = $rEnquete2.Find.Execute(
    "Enquete para saber se os clientes estão satisfeitos com o produto e o serviço que foi prestado.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundEnquete2) {
    $enqueteParaStart = $rEnquete2.Paragraphs(1).Range.Start
    $nextParaIndex = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Start -eq $enqueteParaStart) {
            $nextParaIndex = $i + 1
            break
        }
    }
    $insertPoint2 = $d.Paragraphs($nextParaIndex).Range.End
    $newRange2 = $d.Range($insertPoint2, $insertPoint2)
    $newRange2.InsertXML($emptyParaXml) | Out-Null
}

# ---------------------------------------------------------------------------
# 10) Anúncios em redes sociais (2nd occurrence / "Anúncios, propagandas ...")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Anúncios, propagandas e ofertadas através das redes sociais.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Anúncios, propagandas e ofertadas através de redes sociais.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 11) Aplicativo mobile
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Aplicativo da empresa onde clientes podem baixar para realizar pedidos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Aplicativo da empresa onde clientes podem baixar, para realizar e acompanhar pedidos.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 12) Marketing para capadura: heading capitalization change
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Marketing para capadura:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Marketing para Capadura:",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 13) Anúncios na inovação de capaduras em TCC.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Anúncios na inovação de capaduras em TCC.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Anúncios sobre capadura de TCC para atrair uma clientela maior.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 14) Franquia para expansão do negócio em outras regiões.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Franquia para expansão do negócio em outras regiões.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Franquia para expansão do negócio em outras regiões, visando um lucro maior.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 15) Criação de logotipos para empresas e arte final.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Criação de logotipos para empresas e arte final.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Criação de logotipos para empresas e design de capadura para TCC.",
    2) | Out-Null

Write-Output "done"
